# Weekly Fruta/Hortaliza update: insert 3 new "Pera" price records for
# "Macroferia Regional de Talca" at the top of the existing block (rows
# 307-309), pushing the previous rows 307-389 down to 310-392.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 307 - this shifts the
# existing data (old rows 307-389) down to rows 310-392 and carries the
# date-formatted style from row 306 onto the new D307:D309 cells.
$ws.Rows("307:309").Insert()

# New row 307
$ws.Cells.Item(307, 1).Value  = 5
$ws.Cells.Item(307, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(307, 3).Value  = "Maule"
$ws.Cells.Item(307, 4).Value  = 44508
$ws.Cells.Item(307, 5).Value  = 7
$ws.Cells.Item(307, 6).Value  = "Fruta"
$ws.Cells.Item(307, 7).Value  = 100104
$ws.Cells.Item(307, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(307, 9).Value  = 100104005
$ws.Cells.Item(307, 10).Value = "Pera"
$ws.Cells.Item(307, 11).Value = "Packham's Triumph"
$ws.Cells.Item(307, 12).Value = "Especial"
$ws.Cells.Item(307, 13).Value = 130
$ws.Cells.Item(307, 14).Value = 12000
$ws.Cells.Item(307, 15).Value = 12000
$ws.Cells.Item(307, 16).Value = 12000
$ws.Cells.Item(307, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(307, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(307, 19).Value = 667
$ws.Cells.Item(307, 20).Value = 18

# New row 308
$ws.Cells.Item(308, 1).Value  = 5
$ws.Cells.Item(308, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(308, 3).Value  = "Maule"
$ws.Cells.Item(308, 4).Value  = 44508
$ws.Cells.Item(308, 5).Value  = 7
$ws.Cells.Item(308, 6).Value  = "Fruta"
$ws.Cells.Item(308, 7).Value  = 100104
$ws.Cells.Item(308, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(308, 9).Value  = 100104005
$ws.Cells.Item(308, 10).Value = "Pera"
$ws.Cells.Item(308, 11).Value = "Packham's Triumph"
$ws.Cells.Item(308, 12).Value = "Primera"
$ws.Cells.Item(308, 13).Value = 160
$ws.Cells.Item(308, 14).Value = 10000
$ws.Cells.Item(308, 15).Value = 10000
$ws.Cells.Item(308, 16).Value = 10000
$ws.Cells.Item(308, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(308, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(308, 19).Value = 556
$ws.Cells.Item(308, 20).Value = 18

# New row 309
$ws.Cells.Item(309, 1).Value  = 5
$ws.Cells.Item(309, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(309, 3).Value  = "Maule"
$ws.Cells.Item(309, 4).Value  = 44508
$ws.Cells.Item(309, 5).Value  = 7
$ws.Cells.Item(309, 6).Value  = "Fruta"
$ws.Cells.Item(309, 7).Value  = 100104
$ws.Cells.Item(309, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(309, 9).Value  = 100104005
$ws.Cells.Item(309, 10).Value = "Pera"
$ws.Cells.Item(309, 11).Value = "Packham's Triumph"
$ws.Cells.Item(309, 12).Value = "Segunda"
$ws.Cells.Item(309, 13).Value = 170
$ws.Cells.Item(309, 14).Value = 7000
$ws.Cells.Item(309, 15).Value = 7000
$ws.Cells.Item(309, 16).Value = 7000
$ws.Cells.Item(309, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(309, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(309, 19).Value = 389
$ws.Cells.Item(309, 20).Value = 18
